# Update Streamlit app with latest changes
# Append 3 new data rows (101-103) to the "EDM DATA" worksheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newRows = @(
    @{ Row = 101; Timestamp = "2024-12-19 00:11:35"; B = -0.1191696863941615;  C = -0.001996384113085996; D = 0.009516338747149777 },
    @{ Row = 102; Timestamp = "2024-12-19 00:11:35"; B = -0.1176796625554562;  C = -0.002015617344165995; D = 0.009487886756095161 },
    @{ Row = 103; Timestamp = "2024-12-19 00:11:36"; B = -0.1178976354750921;  C = -0.001868179443807995; D = 0.008810157562725414 }
)

foreach ($r in $newRows) {
    $ws.Cells.Item($r.Row, 1).Value = $r.Timestamp
    $ws.Cells.Item($r.Row, 2).Value = $r.B
    $ws.Cells.Item($r.Row, 3).Value = $r.C
    $ws.Cells.Item($r.Row, 4).Value = $r.D
}
